# ---------------------------------------------------------------------------
# edit.ps1 -- "corrected data cleaning for pre/post/total fixation data"
#
# Summary of the change applied to analysis/pre_gemini_data/participant180/easy.xlsx:
#   1. Header row (row 1) loses the bold/centered/bordered style - cells fall
#      back to the default (unstyled) cell format; A1's leftover
#      "Unnamed: 0" label text is cleared.
#   2. Numeric data in rows 3-7 (Revisit count, Fixation count, Dwell time
#      (ms), Dwell time (%), Fixation duration (ms)) is replaced with
#      recomputed values; the "method" (L), "parameter" (O) and "var3" (U)
#      columns are blanked out for those rows.
#   3. Row 8 (First fixation duration (ms)) keeps its other values but the
#      "method" (L), "parameter" (O) and "var3" (U) columns are blanked too.
#   4. The trailing all-blank rows 10-12 are removed, shrinking the used
#      range from A1:U12 down to A1:U9.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the header styling (bold font / thin border / centered
#        alignment) back to the workbook's default style, and clear the
#        stray "Unnamed: 0" label left in A1. -------------------------------
$ws.Range("A1:U1").ClearFormats()
$ws.Range("A1").ClearContents()

# --- 2. Rewrite rows 3-7 (Revisit count / Fixation count / Dwell time (ms) /
#        Dwell time (%) / Fixation duration (ms)) with the recomputed
#        figures. Columns line up B..U => arg, code, conditional body,
#        conditional statement, external, gemini, literal, literal2,
#        loop body, loop statement, method, method declaration, method2,
#        parameter, return, return2, summary, var, var2, var3. The
#        "method" (L), "parameter" (O) and "var3" (U) columns are blanked.
# Row 3: Revisit count
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 5
$arr[0,1] = 29
$arr[0,2] = 9
$arr[0,3] = 18
$arr[0,4] = 1
$arr[0,5] = 26
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 31
$arr[0,9] = 2
$arr[0,10] = ""
$arr[0,11] = 3
$arr[0,12] = 2
$arr[0,13] = ""
$arr[0,14] = 10
$arr[0,15] = 12
$arr[0,16] = 38
$arr[0,17] = 2
$arr[0,18] = 2
$arr[0,19] = ""
$ws.Range("B3:U3").Value = $arr

# Row 4: Fixation count
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 9
$arr[0,1] = 127
$arr[0,2] = 11
$arr[0,3] = 36
$arr[0,4] = 2
$arr[0,5] = 232
$arr[0,6] = 6
$arr[0,7] = 2
$arr[0,8] = 95
$arr[0,9] = 3
$arr[0,10] = ""
$arr[0,11] = 5
$arr[0,12] = 3
$arr[0,13] = ""
$arr[0,14] = 15
$arr[0,15] = 20
$arr[0,16] = 288
$arr[0,17] = 4
$arr[0,18] = 3
$arr[0,19] = ""
$ws.Range("B4:U4").Value = $arr

# Row 5: Dwell time (ms)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 1835.19
$arr[0,1] = 30660.28
$arr[0,2] = 2769.98
$arr[0,3] = 9226.82
$arr[0,4] = 300.39
$arr[0,5] = 75096.18
$arr[0,6] = 1584.96
$arr[0,7] = 450.35
$arr[0,8] = 23252.82
$arr[0,9] = 583.94
$arr[0,10] = ""
$arr[0,11] = 1201.42
$arr[0,12] = 717.44
$arr[0,13] = ""
$arr[0,14] = 3470.74
$arr[0,15] = 7234.47
$arr[0,16] = 87179.29
$arr[0,17] = 884.19
$arr[0,18] = 583.94
$arr[0,19] = ""
$ws.Range("B5:U5").Value = $arr

# Row 6: Dwell time (%)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.92
$arr[0,1] = 15.4
$arr[0,2] = 1.39
$arr[0,3] = 4.63
$arr[0,4] = 0.15
$arr[0,5] = 40.63
$arr[0,6] = 0.8
$arr[0,7] = 0.23
$arr[0,8] = 11.68
$arr[0,9] = 0.29
$arr[0,10] = ""
$arr[0,11] = 0.6
$arr[0,12] = 0.36
$arr[0,13] = ""
$arr[0,14] = 1.74
$arr[0,15] = 3.63
$arr[0,16] = 43.79
$arr[0,17] = 0.44
$arr[0,18] = 0.29
$arr[0,19] = ""
$ws.Range("B6:U6").Value = $arr

# Row 7: Fixation duration (ms)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 203.91
$arr[0,1] = 241.42
$arr[0,2] = 251.82
$arr[0,3] = 256.3
$arr[0,4] = 150.19
$arr[0,5] = 323.69
$arr[0,6] = 264.16
$arr[0,7] = 225.18
$arr[0,8] = 244.77
$arr[0,9] = 194.65
$arr[0,10] = ""
$arr[0,11] = 240.28
$arr[0,12] = 239.15
$arr[0,13] = ""
$arr[0,14] = 231.38
$arr[0,15] = 361.72
$arr[0,16] = 302.71
$arr[0,17] = 221.05
$arr[0,18] = 194.65
$arr[0,19] = ""
$ws.Range("B7:U7").Value = $arr

# --- 3. Row 8 ("First fixation duration (ms)") - only the method/parameter/
#        var3 columns are cleared; everything else is left as-is. ----------
$ws.Range("L8").ClearContents()
$ws.Range("O8").ClearContents()
$ws.Range("U8").ClearContents()

# --- 4. Drop the trailing blank rows 10-12 so the sheet's used range
#        shrinks from A1:U12 to A1:U9. --------------------------------------
$ws.Range("A10:A12").EntireRow.Delete()
